$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force Text number format on price cells whose new values look like plain
# numbers, so Excel stores them verbatim (matching the source data which uses
# a "1.001"-style string, not a numeric 1.001) instead of silently coercing
# them into Number cells.
$textCells = @(
    "D4",
    "D5",
    "D6",
    "D7",
    "D8",
    "D9",
    "D10",
    "D11",
    "D12",
    "D13",
    "D14",
    "D15",
    "D16",
    "D18",
    "D19",
    "D20",
    "D21",
    "D23",
    "D25",
    "D26",
    "D27",
    "D28",
    "D29",
    "D30",
    "D31",
    "D32",
    "D34",
    "D35",
    "D36",
    "D37",
    "D38",
    "D39",
    "D40",
    "D41",
    "D42",
    "D43",
    "D44",
    "D45",
    "D46",
    "D47",
    "D48",
    "D50",
    "D51"
)
foreach ($cellRef in $textCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

# Apply the updated values.
$ws.Range("D2").Value = '24.437.06'
$ws.Range("E2").Value = '  -1.46%  '
$ws.Range("D3").Value = '1.657.41'
$ws.Range("E3").Value = '  -2.54%  '
$ws.Range("D4").Value = '1.000'
$ws.Range("E4").Value = '  +0.13%  '
$ws.Range("D5").Value = '311.41'
$ws.Range("E5").Value = '  -1.42%  '
$ws.Range("D6").Value = '1.000'
$ws.Range("E6").Value = '  +0.11%  '
$ws.Range("D7").Value = '0.3919'
$ws.Range("E7").Value = '  -1.50%  '
$ws.Range("D8").Value = '0.3925'
$ws.Range("E8").Value = '  -2.59%  '
$ws.Range("D9").Value = '1.000'
$ws.Range("E9").Value = '  +0.23%  '
$ws.Range("D10").Value = '1.385'
$ws.Range("E10").Value = '  -5.90%  '
$ws.Range("D11").Value = '50.38'
$ws.Range("E11").Value = '  -4.55%  '
$ws.Range("D12").Value = '0.08548'
$ws.Range("E12").Value = '  -3.06%  '
$ws.Range("D13").Value = '24.97'
$ws.Range("D14").Value = '7.250'
$ws.Range("E14").Value = '  -3.70%  '
$ws.Range("D15").Value = '0.00001306'
$ws.Range("E15").Value = '  -2.98%  '
$ws.Range("D16").Value = '7.640'
$ws.Range("E16").Value = '  -4.49%  '
$ws.Range("D17").Value = '1.659.61'
$ws.Range("E17").Value = '  -2.38%  '
$ws.Range("D18").Value = '93.72'
$ws.Range("E18").Value = '  -2.16%  '
$ws.Range("D19").Value = '0.06955'
$ws.Range("E19").Value = '  -3.29%  '
$ws.Range("D20").Value = '20.95'
$ws.Range("E20").Value = '  +0.81%  '
$ws.Range("D21").Value = '7.028'
$ws.Range("E21").Value = '  -4.30%  '
$ws.Range("E22").Value = '  +0.11%  '
$ws.Range("D23").Value = '13.86'
$ws.Range("E23").Value = '  -3.61%  '
$ws.Range("D24").Value = '24.434.60'
$ws.Range("E24").Value = '  -1.41%  '
$ws.Range("D25").Value = '2.346'
$ws.Range("E25").Value = '  -0.83%  '
$ws.Range("D26").Value = '2.794'
$ws.Range("E26").Value = '  -4.23%  '
$ws.Range("D27").Value = '22.80'
$ws.Range("E27").Value = '  -1.42%  '
$ws.Range("D28").Value = '159.41'
$ws.Range("E28").Value = '  -1.28%  '
$ws.Range("D29").Value = '5.712'
$ws.Range("E29").Value = '  -7.84%  '
$ws.Range("D30").Value = '145.27'
$ws.Range("E30").Value = '  +0.70%  '
$ws.Range("D31").Value = '8.188'
$ws.Range("E31").Value = '  -4.95%  '
$ws.Range("D32").Value = '2.630'
$ws.Range("E32").Value = '  +8.28%  '
$ws.Range("D33").Value = '1.839.18'
$ws.Range("E33").Value = '  -2.62%  '
$ws.Range("B34").Value = 'ImmutableX'
$ws.Range("C34").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D34").Value = '1.021'
$ws.Range("E34").Value = '  -2.12%  '
$ws.Range("B35").Value = 'Hedera'
$ws.Range("C35").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D35").Value = '0.08249'
$ws.Range("E35").Value = '  -4.68%  '
$ws.Range("D36").Value = '0.03018'
$ws.Range("E36").Value = '  -5.13%  '
$ws.Range("D37").Value = '6.888'
$ws.Range("D38").Value = '0.2783'
$ws.Range("E38").Value = '  -2.16%  '
$ws.Range("D39").Value = '0.09485'
$ws.Range("E39").Value = '  +0.51%  '
$ws.Range("D40").Value = '10.27'
$ws.Range("E40").Value = '  -4.60%  '
$ws.Range("D41").Value = '1.493'
$ws.Range("E41").Value = '  +0.86%  '
$ws.Range("D42").Value = '0.7821'
$ws.Range("E42").Value = '  -5.72%  '
$ws.Range("D43").Value = '13.46'
$ws.Range("E43").Value = '  -5.16%  '
$ws.Range("D44").Value = '16.39'
$ws.Range("E44").Value = '  -6.74%  '
$ws.Range("D45").Value = '2.562'
$ws.Range("E45").Value = '  -4.66%  '
$ws.Range("D46").Value = '0.7049'
$ws.Range("E46").Value = '  -5.01%  '
$ws.Range("D47").Value = '4.154'
$ws.Range("E47").Value = '  -1.44%  '
$ws.Range("D48").Value = '0.08634'
$ws.Range("E48").Value = '  +3.02%  '
$ws.Range("E49").Value = '  +0.17%  '
$ws.Range("D50").Value = '1.313'
$ws.Range("E50").Value = '  -5.54%  '
$ws.Range("D51").Value = '136.85'
$ws.Range("E51").Value = '  -1.99%  '
